# Output file path is removed from Input sheet for NI Scenarios
#
# The "TestResultExcelFilePath" column (column H) is removed from both the
# "ProcessPayrollForNIMonthly" and "TestReports" sheets, shifting the
# columns that followed it one place to the left.

$wb = $excel.ActiveWorkbook

# --- ProcessPayrollForNIMonthly: delete column H (TestResultExcelFilePath) ---
$ws3 = $wb.Worksheets.Item("ProcessPayrollForNIMonthly")
$ws3.Columns.Item(8).Delete()

# --- TestReports: delete column H (TestResultExcelFilePath) ---
$ws4 = $wb.Worksheets.Item("TestReports")
$ws4.Columns.Item(8).Delete()

# The sheet being worked on (ProcessPayrollForNIMonthly) becomes the active tab
$ws4.Range("N2").Select()

$ws3.Activate()
$ws3.Range("M5").Select()
